$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. ALPHASITIO (row 3) Saldo update.
$ws.Range("C3").Value = 300810.78

# 2. Move PHYLIA (row 11) so it comes right before THOMAS (row 10): copy row
#    11 and insert the copy at row 10. This pushes THOMAS down to row 11 and
#    the (now duplicate) original PHYLIA row down to row 12.
$ws.Rows(11).Copy()
$ws.Rows(10).Insert()

# 3. Delete PEDRO / 005324840 row. Originally row 14, shifted to row 15 by
#    the insert above. Delete bottom-most rows first to keep the remaining
#    row numbers below stable.
$ws.Rows(15).Delete()

# 4. Delete the now-duplicate PHYLIA row left behind by the copy (row 12).
$ws.Rows(12).Delete()

# 5. Delete MARIA / 004212581 row (row 6).
$ws.Rows(6).Delete()

# After all the row deletions/inserts above, THIAGO sits at row 7, and
# THOMAS sits at row 10.

# 6. THIAGO Saldo update.
$ws.Range("C7").Value = 24910.96

# 7. THOMAS Saldo update.
$ws.Range("C10").Value = 14567.86
